$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.949.00'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '1.907.03'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8017'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.77'
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +3.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.31'
$ws.Range("E9").Value = '  +5.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06915'
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07987'
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").Value = '1.909.67'
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7416'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.189'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.83'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").Value = '29.957.18'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.95'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.857'
$ws.Range("E18").Value = '  -2.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.13'
$ws.Range("E19").Value = '  +5.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007738'
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '2.153.48'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9992'
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.816'
$ws.Range("E24").Value = '  -1.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.82'
$ws.Range("E25").Value = '  +2.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.206'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1410'
$ws.Range("E27").Value = '  +10.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.89'
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.027'
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.362'
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.514'
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.303'
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.082'
$ws.Range("E33").Value = '  +2.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05508'
$ws.Range("E34").Value = '  +3.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.253'
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7284'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01920'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("E39").Value = '  +0.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.137'
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4414'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.01'
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8314'
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.874'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.57'
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.515'
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.697'
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '990.81'
$ws.Range("E49").Value = '  +8.59%  '
$ws.Range("D50").Value = '2.062.20'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.26'
$ws.Range("E51").Value = '  +0.97%  '
